$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 2
$endRow = 303

for ($r = $startRow; $r -le $endRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value2 = $cell.Value2 + 1
}
